$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy style from an existing header cell (AB1) to the new header cells
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-39: Wins=75, Losses=87, Ties=0
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 29).Value = 75   # AC
    $ws.Cells.Item($r, 30).Value = 87   # AD
    $ws.Cells.Item($r, 31).Value = 0    # AE
}
